$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8257.892
$ws.Range("I132").Value = 7572.7085
$ws.Range("J132").Value = 9522.846
$ws.Range("K132").Value = 22718.1255
$ws.Range("L132").Value = 28568.538
$ws.Range("M132").Value = -20188.1255
$ws.Range("N132").Value = -33628.538

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 25400
$ws.Range("J130").Value = 25400
$ws.Range("L130").Value = 25400
$ws.Range("N130").Value = -35440

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 35956.668
$ws.Range("J132").Value = 35956.668
$ws.Range("L132").Value = 35956.668
$ws.Range("N132").Value = -46076.668

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5053085.5
$ws.Range("I31").Value = 2549.318
$ws.Range("J31").Value = 7578353.5
$ws.Range("K31").Value = 2549.318
$ws.Range("L31").Value = 7578353.5
$ws.Range("M31").Value = -2254.318
$ws.Range("N31").Value = -7578943.5
$ws.Range("H34").Value = 5053085.5
$ws.Range("I34").Value = 2549.318
$ws.Range("J34").Value = 7578353.5
$ws.Range("K34").Value = 2549.318
$ws.Range("L34").Value = 7578353.5
$ws.Range("M34").Value = -2347.318
$ws.Range("N34").Value = -7578757.5
$ws.Range("H122").Value = 3044.25
$ws.Range("I122").Value = 725.6667
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 2177.0001
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = 272.9998999999998
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 2683.625
$ws.Range("I132").Value = 2042.6154
$ws.Range("J132").Value = 3122.2104
$ws.Range("K132").Value = 6127.8462
$ws.Range("L132").Value = 9366.6312
$ws.Range("M132").Value = -3597.8462
$ws.Range("N132").Value = -14426.6312
$ws.Range("H134").Value = 1690.875
$ws.Range("I134").Value = 1280.05
$ws.Range("J134").Value = 2101.7
$ws.Range("K134").Value = 3840.15
$ws.Range("L134").Value = 6305.099999999999
$ws.Range("M134").Value = -1305.15
$ws.Range("N134").Value = -11375.1

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9557.857
$ws.Range("I3").Value = 1962
$ws.Range("J3").Value = 13777.777
$ws.Range("K3").Value = 5886
$ws.Range("L3").Value = 41333.331
$ws.Range("M3").Value = -5774
$ws.Range("N3").Value = -41557.331
$ws.Range("H5").Value = 1140.125
$ws.Range("I5").Value = 333.33334
$ws.Range("J5").Value = 1624.2
$ws.Range("K5").Value = 1000.00002
$ws.Range("L5").Value = 4872.6
$ws.Range("M5").Value = -888.0000200000001
$ws.Range("N5").Value = -5096.6
$ws.Range("H6").Value = 904.8570999999999
$ws.Range("I6").Value = 222.33333
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 666.99999
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = -553.99999
$ws.Range("N6").Value = -15226
$ws.Range("H10").Value = 455.6
$ws.Range("I10").Value = 547.75
$ws.Range("J10").Value = 87
$ws.Range("K10").Value = 1643.25
$ws.Range("L10").Value = 261
$ws.Range("M10").Value = -1504.25
$ws.Range("N10").Value = -539
$ws.Range("H11").Value = 433.25
$ws.Range("I11").Value = 24
$ws.Range("J11").Value = 678.8
$ws.Range("K11").Value = 72
$ws.Range("L11").Value = 2036.4
$ws.Range("M11").Value = 68
$ws.Range("N11").Value = -2316.4
$ws.Range("H13").Value = 1963.8462
$ws.Range("I13").Value = 757.5
$ws.Range("J13").Value = 2500
$ws.Range("K13").Value = 2272.5
$ws.Range("L13").Value = 7500
$ws.Range("M13").Value = -2104.5
$ws.Range("N13").Value = -7836
$ws.Range("H18").Value = 1083.25
$ws.Range("I18").Value = 650
$ws.Range("J18").Value = 1516.5
$ws.Range("K18").Value = 1950
$ws.Range("L18").Value = 4549.5
$ws.Range("M18").Value = -1781
$ws.Range("N18").Value = -4887.5
$ws.Range("H37").Value = 98955.17999999999
$ws.Range("J37").Value = 98955.17999999999
$ws.Range("L37").Value = 296865.54
$ws.Range("N37").Value = -297089.54
$ws.Range("H127").Value = 10101779
$ws.Range("J127").Value = 11364414
$ws.Range("L127").Value = 34093242
$ws.Range("N127").Value = -34103162
$ws.Range("H130").Value = 2308.889
$ws.Range("I130").Value = 2156
$ws.Range("J130").Value = 2500
$ws.Range("K130").Value = 6468
$ws.Range("L130").Value = 7500
$ws.Range("M130").Value = -1448
$ws.Range("N130").Value = -17540
$ws.Range("H131").Value = 910.5848999999999
$ws.Range("I131").Value = 445.29413
$ws.Range("J131").Value = 1130.3055
$ws.Range("K131").Value = 1335.88239
$ws.Range("L131").Value = 3390.9165
$ws.Range("M131").Value = 3704.11761
$ws.Range("N131").Value = -13470.9165
$ws.Range("H132").Value = 959.38464
$ws.Range("J132").Value = 1455
$ws.Range("L132").Value = 13095
$ws.Range("N132").Value = -18155
$ws.Range("H133").Value = 3528.75
$ws.Range("I133").Value = 2007.5
$ws.Range("J133").Value = 5050
$ws.Range("K133").Value = 6022.5
$ws.Range("L133").Value = 15150
$ws.Range("M133").Value = -962.5
$ws.Range("N133").Value = -25270
$ws.Range("H134").Value = 2541.7778
$ws.Range("I134").Value = 1702.5
$ws.Range("J134").Value = 3213.2
$ws.Range("K134").Value = 5107.5
$ws.Range("L134").Value = 9639.599999999999
$ws.Range("M134").Value = -37.5
$ws.Range("N134").Value = -19779.6
$ws.Range("H135").Value = 1140.125
$ws.Range("I135").Value = 333.33334
$ws.Range("J135").Value = 1624.2
$ws.Range("K135").Value = 3000.00006
$ws.Range("L135").Value = 14617.8
$ws.Range("M135").Value = -465.0000600000003
$ws.Range("N135").Value = -19687.8
$ws.Range("H136").Value = 3463.75
$ws.Range("I136").Value = 1927.5
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 5782.5
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -682.5
$ws.Range("N136").Value = -25200
$ws.Range("H137").Value = 4594.1665
$ws.Range("I137").Value = 1466.3636
$ws.Range("J137").Value = 39000
$ws.Range("K137").Value = 4399.0908
$ws.Range("L137").Value = 117000
$ws.Range("M137").Value = 700.9092000000001
$ws.Range("N137").Value = -127200
$ws.Range("H138").Value = 3000.0908
$ws.Range("I138").Value = 1888.6875
$ws.Range("J138").Value = 5963.8335
$ws.Range("K138").Value = 5666.0625
$ws.Range("L138").Value = 17891.5005
$ws.Range("M138").Value = -526.0625
$ws.Range("N138").Value = -28171.5005
$ws.Range("H139").Value = 2461.4092
$ws.Range("I139").Value = 1171.6666
$ws.Range("J139").Value = 2945.0625
$ws.Range("K139").Value = 3514.9998
$ws.Range("L139").Value = 8835.1875
$ws.Range("M139").Value = 1625.0002
$ws.Range("N139").Value = -19115.1875
$ws.Range("H140").Value = 1206.0952
$ws.Range("I140").Value = 1073.7778
$ws.Range("J140").Value = 2000
$ws.Range("K140").Value = 3221.3334
$ws.Range("L140").Value = 6000
$ws.Range("M140").Value = 1958.6666
$ws.Range("N140").Value = -16360
$ws.Range("H141").Value = 2734.8333
$ws.Range("I141").Value = 1222.25
$ws.Range("J141").Value = 5760
$ws.Range("K141").Value = 3666.75
$ws.Range("L141").Value = 17280
$ws.Range("M141").Value = 1513.25
$ws.Range("N141").Value = -27640

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2752.4
$ws.Range("I80").Value = 2764.4443
$ws.Range("J80").Value = 2734.3333
$ws.Range("K80").Value = 2764.4443
$ws.Range("L80").Value = 2734.3333
$ws.Range("M80").Value = -1766.4443
$ws.Range("N80").Value = -4730.3333
$ws.Range("H83").Value = 2752.4
$ws.Range("I83").Value = 2764.4443
$ws.Range("J83").Value = 2734.3333
$ws.Range("K83").Value = 13822.2215
$ws.Range("L83").Value = 13671.6665
$ws.Range("M83").Value = -8830.2215
$ws.Range("N83").Value = -23655.6665

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10599.929
$ws.Range("I122").Value = 13822.667
$ws.Range("J122").Value = 4799
$ws.Range("K122").Value = 41468.001
$ws.Range("L122").Value = 14397
$ws.Range("M122").Value = -39018.001
$ws.Range("N122").Value = -19297

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 25850
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 25850
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 25850
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -35670
$ws.Range("H140").Value = 34440.637
$ws.Range("J140").Value = 34440.637
$ws.Range("L140").Value = 34440.637
$ws.Range("N140").Value = -44800.637
